$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new test-case row (row 6) with its Data payload
$ws.Range("A6").Value = "test_create_user"
$ws.Range("B6").Value = "TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond,`nBusiness1,Business2,+917975935256,+918105855417"
$ws.Range("C6").Value = "Y"

# The Data cell holds a multi-field, multi-line payload - wrap it and size the row to fit
$ws.Range("B6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 58

# Widen column B (drop the old bestFit width) so the wrapped text reads cleanly
$ws.Columns.Item(2).ColumnWidth = 40.5

# Leave the selection where the author's last save left it
$ws.Range("E12").Select() | Out-Null
